$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149. This shifts the existing rows 149-182
# down to 150-183, preserving all their values and formatting (so the old
# row 182 becomes the new row 183 automatically, matching the target file).
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with this week's data (other columns - A, B, C,
# E, F, G, H, I, N, O, Q, R - are identical for every row in this sheet, and
# Excel's Insert already left the row blank, so only the weekly values need
# to be written).
$ws.Range("A149").Value = 8
$ws.Range("B149").Value = "Terminal La Palmera de La Serena"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44889
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100112040
$ws.Range("G149").Value = "Cilantro"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 2400
$ws.Range("K149").Value = 1500
$ws.Range("L149").Value = 2000
$ws.Range("M149").Value = 1750
$ws.Range("N149").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 1167
$ws.Range("Q149").Value = 1.5
$ws.Range("R149").Value = "Hortaliza"

# Make sure the new row's date cell carries the same date format as the
# other rows in column D.
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat()
